# Auto-generated Excel COM-interop script
# Applies the cached-value corrections described in the commit diff
# (re-priced Leve profit calculations across ALC/BSM/CUL/GSM/LTW sheets).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 1814.8
$ws.Range("I2").Value = 1960.8889
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 1960.8889
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = -1847.8889
$ws.Range("N2").Value = -726
# Row 9
$ws.Range("H9").Value = 519.4167
$ws.Range("I9").Value = 651.3214
$ws.Range("J9").Value = 57.75
$ws.Range("K9").Value = 651.3214
$ws.Range("L9").Value = 57.75
$ws.Range("M9").Value = -482.3214
$ws.Range("N9").Value = -395.75
# Row 11
$ws.Range("H11").Value = 65.77778000000001
$ws.Range("I11").Value = 65.77778000000001
$ws.Range("K11").Value = 65.77778000000001
$ws.Range("M11").Value = 74.22221999999999
# Row 12
$ws.Range("H12").Value = 489.2857
$ws.Range("I12").Value = 500.6
$ws.Range("J12").Value = 461
$ws.Range("K12").Value = 500.6
$ws.Range("L12").Value = 461
$ws.Range("M12").Value = -330.6
$ws.Range("N12").Value = -801
# Row 28
$ws.Range("H28").Value = 669.8889
$ws.Range("I28").Value = 691.125
$ws.Range("J28").Value = 500
$ws.Range("K28").Value = 691.125
$ws.Range("L28").Value = 500
$ws.Range("M28").Value = -206.125
$ws.Range("N28").Value = -1470
# Row 29
$ws.Range("H29").Value = 2012.5
$ws.Range("I29").Value = 950
$ws.Range("J29").Value = 3075
$ws.Range("K29").Value = 2850
$ws.Range("L29").Value = 9225
$ws.Range("M29").Value = -2569
$ws.Range("N29").Value = -9787
# Row 33
$ws.Range("H33").Value = 765.2857
$ws.Range("I33").Value = 705.17645
$ws.Range("J33").Value = 1020.75
$ws.Range("K33").Value = 705.17645
$ws.Range("L33").Value = 1020.75
$ws.Range("M33").Value = -476.17645
$ws.Range("N33").Value = -1478.75
# Row 38
$ws.Range("H38").Value = 598.875
$ws.Range("I38").Value = 598.875
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 1796.625
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -1424.625
$ws.Range("N38").ClearContents()
# Row 55
$ws.Range("H55").Value = 239.15384
$ws.Range("I55").Value = 228.77777
$ws.Range("J55").Value = 262.5
$ws.Range("K55").Value = 228.77777
$ws.Range("L55").Value = 262.5
$ws.Range("M55").Value = -14.77777
$ws.Range("N55").Value = -690.5
# Row 87
$ws.Range("H87").Value = 29828.334
$ws.Range("J87").Value = 29828.334
$ws.Range("L87").Value = 29828.334
$ws.Range("N87").Value = -32324.334
# Row 90
$ws.Range("H90").Value = 29828.334
$ws.Range("J90").Value = 29828.334
$ws.Range("L90").Value = 89485.00199999999
$ws.Range("N90").Value = -101965.002
# Row 92
$ws.Range("H92").Value = 6301.6
$ws.Range("I92").Value = 6890.6665
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 6890.6665
$ws.Range("L92").Value = 1000
$ws.Range("M92").Value = -5642.6665
$ws.Range("N92").Value = -3496
# Row 98
$ws.Range("H98").Value = 1383.2703
$ws.Range("I98").Value = 811.6
$ws.Range("K98").Value = 811.6
$ws.Range("M98").Value = 686.4
# Row 107
$ws.Range("H107").Value = 1442
$ws.Range("I107").Value = 2035
$ws.Range("J107").Value = 849
$ws.Range("K107").Value = 2035
$ws.Range("L107").Value = 849
$ws.Range("M107").Value = -115
$ws.Range("N107").Value = -4689
# Row 113
$ws.Range("H113").Value = 4082.2856
$ws.Range("I113").Value = 3726.6667
$ws.Range("J113").Value = 4224.533
$ws.Range("K113").Value = 3726.6667
$ws.Range("L113").Value = 4224.533
$ws.Range("M113").Value = -472.6667000000002
$ws.Range("N113").Value = -10732.533
# Row 122
$ws.Range("H122").Value = 1383.2703
$ws.Range("I122").Value = 811.6
$ws.Range("K122").Value = 2434.8
$ws.Range("M122").Value = 15.19999999999982

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1225.8462
$ws.Range("I20").Value = 958
$ws.Range("J20").Value = 1538.3334
$ws.Range("K20").Value = 958
$ws.Range("L20").Value = 1538.3334
$ws.Range("M20").Value = -711
$ws.Range("N20").Value = -2032.3334
# Row 64
$ws.Range("H64").Value = 450.3846
$ws.Range("I64").Value = 450.85715
$ws.Range("J64").Value = 449.83334
$ws.Range("K64").Value = 450.85715
$ws.Range("L64").Value = 449.83334
$ws.Range("M64").Value = -225.85715
$ws.Range("N64").Value = -899.83334
# Row 67
$ws.Range("H67").Value = 450.3846
$ws.Range("I67").Value = 450.85715
$ws.Range("J67").Value = 449.83334
$ws.Range("K67").Value = 450.85715
$ws.Range("L67").Value = 449.83334
$ws.Range("M67").Value = 329.14285
$ws.Range("N67").Value = -2009.83334
# Row 80
$ws.Range("H80").Value = 108.5
$ws.Range("I80").Value = 118.85714
$ws.Range("J80").Value = 94
$ws.Range("K80").Value = 118.85714
$ws.Range("L80").Value = 94
$ws.Range("M80").Value = 879.14286
$ws.Range("N80").Value = -2090
# Row 83
$ws.Range("H83").Value = 108.5
$ws.Range("I83").Value = 118.85714
$ws.Range("J83").Value = 94
$ws.Range("K83").Value = 594.2857
$ws.Range("L83").Value = 470
$ws.Range("M83").Value = 4397.7143
$ws.Range("N83").Value = -10454

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 54.4375
$ws.Range("I2").Value = 8.875
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 53.25
$ws.Range("L2").Value = 600
$ws.Range("M2").Value = 59.75
$ws.Range("N2").Value = -826
# Row 3
$ws.Range("H3").Value = 4045.9
$ws.Range("I3").Value = 4045.9
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 12137.7
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -12025.7
$ws.Range("N3").ClearContents()
# Row 33
$ws.Range("H33").Value = 196.41667
$ws.Range("I33").Value = 142.88889
$ws.Range("J33").Value = 357
$ws.Range("K33").Value = 857.33334
$ws.Range("L33").Value = 2142
$ws.Range("M33").Value = -574.33334
$ws.Range("N33").Value = -2708
# Row 38
$ws.Range("H38").Value = 106
$ws.Range("I38").Value = 73.333336
$ws.Range("J38").Value = 400
$ws.Range("K38").Value = 220.000008
$ws.Range("L38").Value = 1200
$ws.Range("M38").Value = 126.999992
$ws.Range("N38").Value = -1894
# Row 39
$ws.Range("H39").Value = 3333.6667
$ws.Range("J39").Value = 4071.4285
$ws.Range("L39").Value = 12214.2855
$ws.Range("N39").Value = -12802.2855
# Row 40
$ws.Range("H40").Value = 216.72728
$ws.Range("I40").Value = 124.588234
$ws.Range("J40").Value = 530
$ws.Range("K40").Value = 498.352936
$ws.Range("L40").Value = 2120
$ws.Range("M40").Value = -429.352936
$ws.Range("N40").Value = -2258
# Row 97
$ws.Range("H97").Value = 307.57144
$ws.Range("I97").Value = 157.66667
$ws.Range("J97").Value = 420
$ws.Range("K97").Value = 473.00001
$ws.Range("L97").Value = 1260
$ws.Range("M97").Value = 22.99998999999997
$ws.Range("N97").Value = -2252
# Row 107
$ws.Range("H107").Value = 478.8421
$ws.Range("I107").Value = 122.44444
$ws.Range("J107").Value = 799.6
$ws.Range("K107").Value = 367.33332
$ws.Range("L107").Value = 2398.8
$ws.Range("M107").Value = 1552.66668
$ws.Range("N107").Value = -6238.8

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 4271.7173
$ws.Range("I80").Value = 4699.971
$ws.Range("K80").Value = 4699.971
$ws.Range("M80").Value = -3701.971
# Row 83
$ws.Range("H83").Value = 4271.7173
$ws.Range("I83").Value = 4699.971
$ws.Range("K83").Value = 23499.855
$ws.Range("M83").Value = -18507.855
# Row 107
$ws.Range("H107").Value = 6894
$ws.Range("I107").Value = 10114
$ws.Range("J107").Value = 454
$ws.Range("K107").Value = 10114
$ws.Range("L107").Value = 454
$ws.Range("M107").Value = -8194
$ws.Range("N107").Value = -4294

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3874.1052
$ws.Range("I40").Value = 3454.818
$ws.Range("J40").Value = 4450.625
$ws.Range("K40").Value = 3454.818
$ws.Range("L40").Value = 4450.625
$ws.Range("M40").Value = -3318.818
$ws.Range("N40").Value = -4722.625

Write-Host "Applied Leve profit cache updates to ALC, BSM, CUL, GSM, LTW sheets."
